$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing text storage (the sheet
# stores Price/Volume figures as text, e.g. "212.27" / "0.130", not numbers -
# several look like plain numbers and would otherwise get auto-converted by
# Excel's normal typed-input parsing). Temporarily mark the cell as Text,
# assign, then restore the default "Normal" style so no stray formatting is
# left behind on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Simple Price (D) / Volume(1h) (E) updates, keyed by row number.
$updates = @{
    2  = @{ D = "88.838.27";  E = "  -1.87%  " }
    3  = @{ D = "3.091.26";   E = "  -3.47%  " }
    4  = @{ E = "  -0.07%  " }
    5  = @{ D = "212.27";     E = "  -4.12%  " }
    6  = @{ D = "621.68";     E = "  -3.20%  " }
    7  = @{ D = "0.374";      E = "  -7.12%  " }
    8  = @{ D = "0.803";      E = "  +13.43%  " }
    9  = @{ E = "  +0.00%  " }
    10 = @{ D = "3.087.54";   E = "  -3.51%  " }
    11 = @{ E = "  +2.54%  " }
    12 = @{ D = "0.181";      E = "  +0.33%  " }
    13 = @{ D = "0.0000242";  E = "  -7.04%  " }
    14 = @{ D = "5.29";       E = "  -2.69%  " }
    15 = @{ D = "88.371.66";  E = "  -2.19%  " }
    16 = @{ D = "32.28";      E = "  -3.72%  " }
    17 = @{ D = "3.662.60";   E = "  -3.49%  " }
    18 = @{ D = "3.071.37";   E = "  -4.28%  " }
    19 = @{ D = "3.38";       E = "  -0.17%  " }
    20 = @{ D = "0.0000209";  E = "  -8.25%  " }
    21 = @{ D = "13.44";      E = "  -0.66%  " }
    22 = @{ D = "423.14";     E = "  -4.06%  " }
    23 = @{ D = "8.26";       E = "  -4.96%  " }
    24 = @{ D = "4.93";       E = "  -3.07%  " }
    25 = @{ D = "5.64";       E = "  +5.40%  " }
    26 = @{ D = "11.93";      E = "  -0.48%  " }
    27 = @{ D = "82.74";      E = "  +1.33%  " }
    28 = @{ D = "3.231.74";   E = "  -4.37%  " }
    29 = @{ D = "1.01";       E = "  +0.74%  " }
    30 = @{ D = "0.171";      E = "  +7.00%  " }
    31 = @{ E = "  +7.90%  " }
    32 = @{ D = "8.12";       E = "  -4.19%  " }
    33 = @{ D = "509.99";     E = "  -6.17%  " }
    34 = @{ D = "3.67";       E = "  -13.82%  " }
    35 = @{ D = "6.80";       E = "  -4.62%  " }
    36 = @{ D = "1.26";       E = "  -3.42%  " }
    37 = @{ D = "1.80";       E = "  -6.40%  " }
    38 = @{ D = "22.28";      E = "  -1.35%  " }
    41 = @{ E = "  +0.13%  " }
    42 = @{ E = "  +0.02%  " }
    43 = @{ D = "0.364";      E = "  -2.97%  " }
    46 = @{ E = "  +4.37%  " }
    51 = @{ D = "0.706";      E = "  -6.36%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        Set-TextValue $ws.Range("D$row") $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        Set-TextValue $ws.Range("E$row") $vals["E"]
    }
}

# Rows whose Coin/Link/Price/Volume data were swapped with an adjacent row.
$rowSwaps = @{
    39 = @{ B = "WhiteBITCoin"; C = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"; D = "22.27";  E = "  -0.62%  " }
    40 = @{ B = "Kaspa";        C = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas";        D = "0.130"; E = "  +1.70%  " }

    44 = @{ B = "Stacks"; C = "https://coinranking.com/coin/mMPrMcB7+stacks-stx";          D = "1.83";   E = "  -6.40%  " }
    45 = @{ B = "Monero"; C = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr";      D = "147.29"; E = "  +0.30%  " }

    47 = @{ B = "Hedera"; C = "https://coinranking.com/coin/jad286TjB+hedera-hbar"; D = "0.0698"; E = "  +14.40%  " }
    48 = @{ B = "OKB";    C = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D = "43.46";  E = "  -3.14%  " }

    49 = @{ B = "Aave";        C = "https://coinranking.com/coin/ixgUfzmLR+aave-aave";      D = "162.42"; E = "  -6.58%  " }
    50 = @{ B = "ImmutableX";  C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx";  D = "1.22";   E = "  -1.94%  " }
}

foreach ($row in $rowSwaps.Keys) {
    $vals = $rowSwaps[$row]
    $ws.Range("B$row").Value = $vals["B"]
    $ws.Range("C$row").Value = $vals["C"]
    Set-TextValue $ws.Range("D$row") $vals["D"]
    Set-TextValue $ws.Range("E$row") $vals["E"]
}
